# Apply updated "dSF" (column F) values for specific rows on Sheet1,
# reflecting a repull/recalculation of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 0
    6  = 1
    10 = -2
    12 = 8
    15 = -3
    16 = -4
    17 = 0
    25 = -14
    37 = 3
    38 = 3
    52 = 1
    53 = -1
    58 = 0
    67 = 2
    70 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
